# This script reproduces a scheduled-runner data refresh for the Omega_Profits
# workbook: it (1) strips the bold/border/center-aligned formatting that used to
# be applied to the header row of every sheet, and (2) refreshes the computed
# market-price/profit columns (H:N) with newly pulled values on the rows where
# prices moved since the last run.

$wb = $excel.ActiveWorkbook

# ---- 1. Clear header-row (row 1) formatting on every sheet ----
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("A1:N1").ClearFormats()
}

# ---- 2. Refresh market data values ----

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(19, 8).Value = 1522.6207  # H19: 1523.5172 -> 1522.6207
$ws.Cells.Item(19, 9).Value = 1243.9375  # I19: 1279.9333 -> 1243.9375
$ws.Cells.Item(19, 10).Value = 1865.6154  # J19: 1784.5 -> 1865.6154
$ws.Cells.Item(19, 11).Value = 1243.9375  # K19: 1279.9333 -> 1243.9375
$ws.Cells.Item(19, 12).Value = 1865.6154  # L19: 1784.5 -> 1865.6154
$ws.Cells.Item(19, 13).Value = -1068.9375  # M19: -1104.9333 -> -1068.9375
$ws.Cells.Item(19, 14).Value = -2215.6154  # N19: -2134.5 -> -2215.6154
$ws.Cells.Item(43, 8).Value = 11448.625  # H43: 11903.158 -> 11448.625
$ws.Cells.Item(43, 9).Value = 8249.75  # I43: 9498 -> 8249.75
$ws.Cells.Item(43, 10).Value = 12514.917  # J43: 12186.117 -> 12514.917
$ws.Cells.Item(43, 11).Value = 8249.75  # K43: 9498 -> 8249.75
$ws.Cells.Item(43, 12).Value = 12514.917  # L43: 12186.117 -> 12514.917
$ws.Cells.Item(43, 13).Value = -8180.75  # M43: -9429 -> -8180.75
$ws.Cells.Item(43, 14).Value = -12652.917  # N43: -12324.117 -> -12652.917
$ws.Cells.Item(57, 8).Value = 0  # H57: 15000 -> 0
$ws.Cells.Item(57, 9).Value = 0  # I57: 15000 -> 0
$ws.Cells.Item(57, 10).Value = 0  # J57: 0 -> 0
$ws.Cells.Item(57, 11).Value = 0  # K57: 45000 -> 0
$ws.Cells.Item(57, 12).Value = 0  # L57: 0 -> 0
$ws.Cells.Item(57, 13).ClearContents()  # M57: -44501 -> (cleared)
$ws.Cells.Item(82, 8).Value = 1291.2858  # H82: 1307.8 -> 1291.2858
$ws.Cells.Item(82, 9).Value = 1291.2858  # I82: 1307.8 -> 1291.2858
$ws.Cells.Item(82, 10).Value = 0  # J82: 0 -> 0
$ws.Cells.Item(82, 11).Value = 3873.8574  # K82: 3923.4 -> 3873.8574
$ws.Cells.Item(82, 12).Value = 0  # L82: 0 -> 0
$ws.Cells.Item(82, 13).Value = -3467.8574  # M82: -3517.4 -> -3467.8574
$ws.Cells.Item(85, 8).Value = 1291.2858  # H85: 1307.8 -> 1291.2858
$ws.Cells.Item(85, 9).Value = 1291.2858  # I85: 1307.8 -> 1291.2858
$ws.Cells.Item(85, 10).Value = 0  # J85: 0 -> 0
$ws.Cells.Item(85, 11).Value = 3873.8574  # K85: 3923.4 -> 3873.8574
$ws.Cells.Item(85, 12).Value = 0  # L85: 0 -> 0
$ws.Cells.Item(85, 13).Value = -2469.8574  # M85: -2519.4 -> -2469.8574
$ws.Cells.Item(113, 8).Value = 3973.5  # H113: 4178.6 -> 3973.5
$ws.Cells.Item(113, 9).Value = 3973.5  # I113: 4178.6 -> 3973.5
$ws.Cells.Item(113, 10).Value = 0  # J113: 0 -> 0
$ws.Cells.Item(113, 11).Value = 3973.5  # K113: 4178.6 -> 3973.5
$ws.Cells.Item(113, 12).Value = 0  # L113: 0 -> 0
$ws.Cells.Item(113, 13).Value = -719.5  # M113: -924.6000000000004 -> -719.5
$ws.Cells.Item(132, 8).Value = 2370.0728  # H132: 2397.7222 -> 2370.0728
$ws.Cells.Item(132, 9).Value = 2234.8235  # I132: 2261.98 -> 2234.8235
$ws.Cells.Item(132, 11).Value = 6704.470499999999  # K132: 6785.940000000001 -> 6704.470499999999
$ws.Cells.Item(132, 13).Value = -4174.470499999999  # M132: -4255.940000000001 -> -4174.470499999999
$ws.Cells.Item(138, 8).Value = 4082.3225  # H138: 4052.4285 -> 4082.3225
$ws.Cells.Item(138, 10).Value = 6341.4243  # J138: 6219.5884 -> 6341.4243
$ws.Cells.Item(138, 12).Value = 19024.2729  # L138: 18658.7652 -> 19024.2729
$ws.Cells.Item(138, 14).Value = -29304.2729  # N138: -28938.7652 -> -29304.2729

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 12351.134  # H32: 12377.167 -> 12351.134
$ws.Cells.Item(32, 9).Value = 1446  # I32: 1518.1765 -> 1446
$ws.Cells.Item(32, 10).Value = 28708.834  # J32: 26577.385 -> 28708.834
$ws.Cells.Item(32, 11).Value = 1446  # K32: 1518.1765 -> 1446
$ws.Cells.Item(32, 12).Value = 28708.834  # L32: 26577.385 -> 28708.834
$ws.Cells.Item(32, 13).Value = -1159  # M32: -1231.1765 -> -1159
$ws.Cells.Item(32, 14).Value = -29282.834  # N32: -27151.385 -> -29282.834
$ws.Cells.Item(61, 8).Value = 4002.6584  # H61: 3866.6 -> 4002.6584
$ws.Cells.Item(61, 9).Value = 3503.6487  # I61: 3295.575 -> 3503.6487
$ws.Cells.Item(61, 10).Value = 8618.5  # J61: 8434.8 -> 8618.5
$ws.Cells.Item(61, 11).Value = 3503.6487  # K61: 3295.575 -> 3503.6487
$ws.Cells.Item(61, 12).Value = 8618.5  # L61: 8434.8 -> 8618.5
$ws.Cells.Item(61, 13).Value = -3291.6487  # M61: -3083.575 -> -3291.6487
$ws.Cells.Item(61, 14).Value = -9042.5  # N61: -8858.8 -> -9042.5
$ws.Cells.Item(74, 8).Value = 2190.3225  # H74: 2334.9429 -> 2190.3225
$ws.Cells.Item(74, 10).Value = 5789  # J74: 3922.4 -> 5789
$ws.Cells.Item(74, 12).Value = 5789  # L74: 3922.4 -> 5789
$ws.Cells.Item(74, 14).Value = -7537  # N74: -5670.4 -> -7537
$ws.Cells.Item(77, 8).Value = 2190.3225  # H77: 2334.9429 -> 2190.3225
$ws.Cells.Item(77, 10).Value = 5789  # J77: 3922.4 -> 5789
$ws.Cells.Item(77, 12).Value = 28945  # L77: 19612 -> 28945
$ws.Cells.Item(77, 14).Value = -37681  # N77: -28348 -> -37681
$ws.Cells.Item(97, 8).Value = 1642.8158  # H97: 1603.1282 -> 1642.8158
$ws.Cells.Item(97, 9).Value = 1063.9642  # I97: 1030.5518 -> 1063.9642
$ws.Cells.Item(97, 11).Value = 1063.9642  # K97: 1030.5518 -> 1063.9642
$ws.Cells.Item(97, 13).Value = -567.9641999999999  # M97: -534.5518 -> -567.9641999999999
$ws.Cells.Item(110, 8).Value = 1549.2727  # H110: 1496.0834 -> 1549.2727
$ws.Cells.Item(110, 9).Value = 1204.2  # I110: 1177.5454 -> 1204.2
$ws.Cells.Item(110, 11).Value = 1204.2  # K110: 1177.5454 -> 1204.2
$ws.Cells.Item(110, 13).Value = 840.8  # M110: 867.4546 -> 840.8
$ws.Cells.Item(132, 8).Value = 4145.4243  # H132: 4076.2646 -> 4145.4243
$ws.Cells.Item(132, 9).Value = 4145.4243  # I132: 4076.2646 -> 4145.4243
$ws.Cells.Item(132, 11).Value = 12436.2729  # K132: 12228.7938 -> 12436.2729
$ws.Cells.Item(132, 13).Value = -9906.2729  # M132: -9698.7938 -> -9906.2729
$ws.Cells.Item(136, 8).Value = 4002.6584  # H136: 3866.6 -> 4002.6584
$ws.Cells.Item(136, 9).Value = 3503.6487  # I136: 3295.575 -> 3503.6487
$ws.Cells.Item(136, 10).Value = 8618.5  # J136: 8434.8 -> 8618.5
$ws.Cells.Item(136, 11).Value = 10510.9461  # K136: 9886.724999999999 -> 10510.9461
$ws.Cells.Item(136, 12).Value = 25855.5  # L136: 25304.4 -> 25855.5
$ws.Cells.Item(136, 13).Value = -7960.946100000001  # M136: -7336.724999999999 -> -7960.946100000001
$ws.Cells.Item(136, 14).Value = -30955.5  # N136: -30404.4 -> -30955.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(26, 8).Value = 0  # H26: 50000 -> 0
$ws.Cells.Item(26, 10).Value = 0  # J26: 50000 -> 0
$ws.Cells.Item(26, 12).Value = 0  # L26: 50000 -> 0
$ws.Cells.Item(26, 14).ClearContents()  # N26: -50584 -> (cleared)
$ws.Cells.Item(44, 8).Value = 0  # H44: 27818 -> 0
$ws.Cells.Item(44, 9).Value = 0  # I44: 18000 -> 0
$ws.Cells.Item(44, 10).Value = 0  # J44: 29999.777 -> 0
$ws.Cells.Item(44, 11).Value = 0  # K44: 18000 -> 0
$ws.Cells.Item(44, 12).Value = 0  # L44: 29999.777 -> 0
$ws.Cells.Item(44, 13).ClearContents()  # M44: -17503 -> (cleared)
$ws.Cells.Item(44, 14).ClearContents()  # N44: -30993.777 -> (cleared)
$ws.Cells.Item(57, 8).Value = 0  # H57: 89996 -> 0
$ws.Cells.Item(57, 10).Value = 0  # J57: 89996 -> 0
$ws.Cells.Item(57, 12).Value = 0  # L57: 89996 -> 0
$ws.Cells.Item(57, 14).ClearContents()  # N57: -91436 -> (cleared)
$ws.Cells.Item(82, 8).Value = 73744.5  # H82: 99992.336 -> 73744.5
$ws.Cells.Item(82, 9).Value = 10000  # I82: 0 -> 10000
$ws.Cells.Item(82, 10).Value = 94992.664  # J82: 99992.336 -> 94992.664
$ws.Cells.Item(82, 11).Value = 10000  # K82: 0 -> 10000
$ws.Cells.Item(82, 12).Value = 94992.664  # L82: 99992.336 -> 94992.664
$ws.Cells.Item(82, 13).Value = -9617  # M82: None -> -9617
$ws.Cells.Item(82, 14).Value = -95758.664  # N82: -100758.336 -> -95758.664
$ws.Cells.Item(85, 8).Value = 73744.5  # H85: 99992.336 -> 73744.5
$ws.Cells.Item(85, 9).Value = 10000  # I85: 0 -> 10000
$ws.Cells.Item(85, 10).Value = 94992.664  # J85: 99992.336 -> 94992.664
$ws.Cells.Item(85, 11).Value = 10000  # K85: 0 -> 10000
$ws.Cells.Item(85, 12).Value = 94992.664  # L85: 99992.336 -> 94992.664
$ws.Cells.Item(85, 13).Value = -8674  # M85: None -> -8674
$ws.Cells.Item(85, 14).Value = -97644.664  # N85: -102644.336 -> -97644.664
$ws.Cells.Item(86, 8).Value = 2922.647  # H86: 3053 -> 2922.647
$ws.Cells.Item(86, 9).Value = 3219.8  # I86: 3182.182 -> 3219.8
$ws.Cells.Item(86, 10).Value = 2498.1428  # J86: 2816.1667 -> 2498.1428
$ws.Cells.Item(86, 11).Value = 3219.8  # K86: 3182.182 -> 3219.8
$ws.Cells.Item(86, 12).Value = 2498.1428  # L86: 2816.1667 -> 2498.1428
$ws.Cells.Item(86, 13).Value = -2096.8  # M86: -2059.182 -> -2096.8
$ws.Cells.Item(86, 14).Value = -4744.1428  # N86: -5062.1667 -> -4744.1428
$ws.Cells.Item(89, 8).Value = 2922.647  # H89: 3053 -> 2922.647
$ws.Cells.Item(89, 9).Value = 3219.8  # I89: 3182.182 -> 3219.8
$ws.Cells.Item(89, 10).Value = 2498.1428  # J89: 2816.1667 -> 2498.1428
$ws.Cells.Item(89, 11).Value = 16099  # K89: 15910.91 -> 16099
$ws.Cells.Item(89, 12).Value = 12490.714  # L89: 14080.8335 -> 12490.714
$ws.Cells.Item(89, 13).Value = -10483  # M89: -10294.91 -> -10483
$ws.Cells.Item(89, 14).Value = -23722.714  # N89: -25312.8335 -> -23722.714
$ws.Cells.Item(94, 8).Value = 1110.4615  # H94: 1118.2307 -> 1110.4615
$ws.Cells.Item(94, 9).Value = 968.36365  # I94: 977.5455 -> 968.36365
$ws.Cells.Item(94, 11).Value = 968.36365  # K94: 977.5455 -> 968.36365
$ws.Cells.Item(94, 13).Value = -517.36365  # M94: -526.5455 -> -517.36365
$ws.Cells.Item(96, 8).Value = 18904.334  # H96: 23232.334 -> 18904.334
$ws.Cells.Item(96, 9).Value = 18904.334  # I96: 23232.334 -> 18904.334
$ws.Cells.Item(96, 11).Value = 18904.334  # K96: 23232.334 -> 18904.334
$ws.Cells.Item(96, 13).Value = -16158.334  # M96: -20486.334 -> -16158.334
$ws.Cells.Item(132, 8).Value = 0  # H132: 60000 -> 0
$ws.Cells.Item(132, 10).Value = 0  # J132: 60000 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 60000 -> 0
$ws.Cells.Item(132, 14).ClearContents()  # N132: -70120 -> (cleared)
$ws.Cells.Item(134, 8).Value = 1466.75  # H134: 1545.6538 -> 1466.75
$ws.Cells.Item(134, 9).Value = 1466.75  # I134: 1545.6538 -> 1466.75
$ws.Cells.Item(134, 11).Value = 4400.25  # K134: 4636.9614 -> 4400.25
$ws.Cells.Item(134, 13).Value = -1865.25  # M134: -2101.9614 -> -1865.25
$ws.Cells.Item(135, 8).Value = 74999  # H135: 72498.75 -> 74999
$ws.Cells.Item(135, 10).Value = 74999  # J135: 72498.75 -> 74999
$ws.Cells.Item(135, 12).Value = 74999  # L135: 72498.75 -> 74999
$ws.Cells.Item(135, 14).Value = -85139  # N135: -82638.75 -> -85139
$ws.Cells.Item(136, 8).Value = 0  # H136: 89996 -> 0
$ws.Cells.Item(136, 10).Value = 0  # J136: 89996 -> 0
$ws.Cells.Item(136, 12).Value = 0  # L136: 89996 -> 0
$ws.Cells.Item(136, 14).ClearContents()  # N136: -100196 -> (cleared)
$ws.Cells.Item(137, 8).Value = 0  # H137: 79710 -> 0
$ws.Cells.Item(137, 10).Value = 0  # J137: 79710 -> 0
$ws.Cells.Item(137, 12).Value = 0  # L137: 79710 -> 0
$ws.Cells.Item(137, 14).ClearContents()  # N137: -89910 -> (cleared)
$ws.Cells.Item(138, 8).Value = 84603.664  # H138: 84581.625 -> 84603.664
$ws.Cells.Item(138, 10).Value = 84603.664  # J138: 84581.625 -> 84603.664
$ws.Cells.Item(138, 12).Value = 84603.664  # L138: 84581.625 -> 84603.664
$ws.Cells.Item(138, 14).Value = -94883.664  # N138: -94861.625 -> -94883.664
$ws.Cells.Item(139, 8).Value = 80709  # H139: 72854 -> 80709
$ws.Cells.Item(139, 10).Value = 0  # J139: 64999 -> 0
$ws.Cells.Item(139, 12).Value = 0  # L139: 64999 -> 0
$ws.Cells.Item(139, 14).ClearContents()  # N139: -75279 -> (cleared)
$ws.Cells.Item(140, 8).Value = 81127.6  # H140: 80616.875 -> 81127.6
$ws.Cells.Item(140, 10).Value = 81127.6  # J140: 80616.875 -> 81127.6
$ws.Cells.Item(140, 12).Value = 81127.6  # L140: 80616.875 -> 81127.6
$ws.Cells.Item(140, 14).Value = -91487.6  # N140: -90976.875 -> -91487.6
$ws.Cells.Item(141, 8).Value = 77499.5  # H141: 79999 -> 77499.5
$ws.Cells.Item(141, 10).Value = 77499.5  # J141: 79999 -> 77499.5
$ws.Cells.Item(141, 12).Value = 77499.5  # L141: 79999 -> 77499.5
$ws.Cells.Item(141, 14).Value = -87859.5  # N141: -90359 -> -87859.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 3277.5476  # H31: 3445.9268 -> 3277.5476
$ws.Cells.Item(31, 9).Value = 2180.2593  # I31: 2327.2 -> 2180.2593
$ws.Cells.Item(31, 10).Value = 5252.6665  # J31: 5193.9375 -> 5252.6665
$ws.Cells.Item(31, 11).Value = 2180.2593  # K31: 2327.2 -> 2180.2593
$ws.Cells.Item(31, 12).Value = 5252.6665  # L31: 5193.9375 -> 5252.6665
$ws.Cells.Item(31, 13).Value = -1885.2593  # M31: -2032.2 -> -1885.2593
$ws.Cells.Item(31, 14).Value = -5842.6665  # N31: -5783.9375 -> -5842.6665
$ws.Cells.Item(34, 8).Value = 3277.5476  # H34: 3445.9268 -> 3277.5476
$ws.Cells.Item(34, 9).Value = 2180.2593  # I34: 2327.2 -> 2180.2593
$ws.Cells.Item(34, 10).Value = 5252.6665  # J34: 5193.9375 -> 5252.6665
$ws.Cells.Item(34, 11).Value = 2180.2593  # K34: 2327.2 -> 2180.2593
$ws.Cells.Item(34, 12).Value = 5252.6665  # L34: 5193.9375 -> 5252.6665
$ws.Cells.Item(34, 13).Value = -1978.2593  # M34: -2125.2 -> -1978.2593
$ws.Cells.Item(34, 14).Value = -5656.6665  # N34: -5597.9375 -> -5656.6665
$ws.Cells.Item(99, 8).Value = 2751.4211  # H99: 2694.3 -> 2751.4211
$ws.Cells.Item(99, 9).Value = 2056.923  # I99: 2096.0833 -> 2056.923
$ws.Cells.Item(99, 10).Value = 4256.1665  # J99: 3591.625 -> 4256.1665
$ws.Cells.Item(99, 11).Value = 2056.923  # K99: 2096.0833 -> 2056.923
$ws.Cells.Item(99, 12).Value = 4256.1665  # L99: 3591.625 -> 4256.1665
$ws.Cells.Item(99, 13).Value = -558.9229999999998  # M99: -598.0832999999998 -> -558.9229999999998
$ws.Cells.Item(99, 14).Value = -7252.1665  # N99: -6587.625 -> -7252.1665
$ws.Cells.Item(126, 8).Value = 2751.4211  # H126: 2694.3 -> 2751.4211
$ws.Cells.Item(126, 9).Value = 2056.923  # I126: 2096.0833 -> 2056.923
$ws.Cells.Item(126, 10).Value = 4256.1665  # J126: 3591.625 -> 4256.1665
$ws.Cells.Item(126, 11).Value = 6170.768999999999  # K126: 6288.249899999999 -> 6170.768999999999
$ws.Cells.Item(126, 12).Value = 12768.4995  # L126: 10774.875 -> 12768.4995
$ws.Cells.Item(126, 13).Value = -3700.768999999999  # M126: -3818.249899999999 -> -3700.768999999999
$ws.Cells.Item(126, 14).Value = -17708.4995  # N126: -15714.875 -> -17708.4995
$ws.Cells.Item(132, 8).Value = 2955.8965  # H132: 2930.6667 -> 2955.8965
$ws.Cells.Item(132, 9).Value = 2511.889  # I132: 2523.037 -> 2511.889
$ws.Cells.Item(132, 10).Value = 8950  # J132: 6599.3335 -> 8950
$ws.Cells.Item(132, 11).Value = 7535.667  # K132: 7569.110999999999 -> 7535.667
$ws.Cells.Item(132, 12).Value = 26850  # L132: 19798.0005 -> 26850
$ws.Cells.Item(132, 13).Value = -5005.667  # M132: -5039.110999999999 -> -5005.667
$ws.Cells.Item(132, 14).Value = -31910  # N132: -24858.0005 -> -31910
$ws.Cells.Item(134, 8).Value = 2715  # H134: 2715.1667 -> 2715
$ws.Cells.Item(134, 9).Value = 2234.5454  # I134: 2234.7273 -> 2234.5454
$ws.Cells.Item(134, 11).Value = 6703.6362  # K134: 6704.1819 -> 6703.6362
$ws.Cells.Item(134, 13).Value = -4168.6362  # M134: -4169.1819 -> -4168.6362
$ws.Cells.Item(135, 8).Value = 0  # H135: 99988 -> 0
$ws.Cells.Item(135, 10).Value = 0  # J135: 99988 -> 0
$ws.Cells.Item(135, 12).Value = 0  # L135: 99988 -> 0
$ws.Cells.Item(135, 14).ClearContents()  # N135: -110128 -> (cleared)

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 8).Value = 4077.75  # H2: 2738.3333 -> 4077.75
$ws.Cells.Item(2, 9).Value = 5066.6665  # I2: 3812.5 -> 5066.6665
$ws.Cells.Item(2, 10).Value = 1111  # J2: 590 -> 1111
$ws.Cells.Item(2, 11).Value = 30399.999  # K2: 22875 -> 30399.999
$ws.Cells.Item(2, 12).Value = 6666  # L2: 3540 -> 6666
$ws.Cells.Item(2, 13).Value = -30286.999  # M2: -22762 -> -30286.999
$ws.Cells.Item(2, 14).Value = -6892  # N2: -3766 -> -6892
$ws.Cells.Item(6, 8).Value = 281.2  # H6: 128.36363 -> 281.2
$ws.Cells.Item(6, 9).Value = 281.2  # I6: 128.36363 -> 281.2
$ws.Cells.Item(6, 11).Value = 843.5999999999999  # K6: 385.09089 -> 843.5999999999999
$ws.Cells.Item(6, 13).Value = -730.5999999999999  # M6: -272.09089 -> -730.5999999999999
$ws.Cells.Item(17, 8).Value = 254.1  # H17: 120.5 -> 254.1
$ws.Cells.Item(17, 9).Value = 116.375  # I17: 122 -> 116.375
$ws.Cells.Item(17, 10).Value = 805  # J17: 110 -> 805
$ws.Cells.Item(17, 11).Value = 349.125  # K17: 366 -> 349.125
$ws.Cells.Item(17, 12).Value = 2415  # L17: 330 -> 2415
$ws.Cells.Item(17, 13).Value = -180.125  # M17: -197 -> -180.125
$ws.Cells.Item(17, 14).Value = -2753  # N17: -668 -> -2753
$ws.Cells.Item(33, 8).Value = 33.875  # H33: 30.222221 -> 33.875
$ws.Cells.Item(33, 9).Value = 33.875  # I33: 30.222221 -> 33.875
$ws.Cells.Item(33, 11).Value = 203.25  # K33: 181.333326 -> 203.25
$ws.Cells.Item(33, 13).Value = 79.75  # M33: 101.666674 -> 79.75
$ws.Cells.Item(39, 8).Value = 2739.5789  # H39: 2746.2273 -> 2739.5789
$ws.Cells.Item(39, 10).Value = 2870.6667  # J39: 2856.9443 -> 2870.6667
$ws.Cells.Item(39, 12).Value = 8612.000100000001  # L39: 8570.832900000001 -> 8612.000100000001
$ws.Cells.Item(39, 14).Value = -9200.000100000001  # N39: -9158.832900000001 -> -9200.000100000001
$ws.Cells.Item(55, 8).Value = 3478  # H55: 3352.5715 -> 3478
$ws.Cells.Item(55, 9).Value = 981  # I55: 869.4286 -> 981
$ws.Cells.Item(55, 10).Value = 5975  # J55: 5835.7144 -> 5975
$ws.Cells.Item(55, 11).Value = 2943  # K55: 2608.2858 -> 2943
$ws.Cells.Item(55, 12).Value = 17925  # L55: 17507.1432 -> 17925
$ws.Cells.Item(55, 13).Value = -2766  # M55: -2431.2858 -> -2766
$ws.Cells.Item(55, 14).Value = -18279  # N55: -17861.1432 -> -18279
$ws.Cells.Item(62, 8).Value = 5250  # H62: 0 -> 5250
$ws.Cells.Item(62, 9).Value = 3000  # I62: 0 -> 3000
$ws.Cells.Item(62, 10).Value = 7500  # J62: 0 -> 7500
$ws.Cells.Item(62, 11).Value = 9000  # K62: 0 -> 9000
$ws.Cells.Item(62, 12).Value = 22500  # L62: 0 -> 22500
$ws.Cells.Item(62, 13).Value = -8314  # M62: None -> -8314
$ws.Cells.Item(62, 14).Value = -23872  # N62: None -> -23872
$ws.Cells.Item(63, 8).Value = 3725  # H63: 0 -> 3725
$ws.Cells.Item(63, 9).Value = 4000  # I63: 0 -> 4000
$ws.Cells.Item(63, 10).Value = 3450  # J63: 0 -> 3450
$ws.Cells.Item(63, 11).Value = 12000  # K63: 0 -> 12000
$ws.Cells.Item(63, 12).Value = 10350  # L63: 0 -> 10350
$ws.Cells.Item(63, 13).Value = -11251  # M63: None -> -11251
$ws.Cells.Item(63, 14).Value = -11848  # N63: None -> -11848
$ws.Cells.Item(64, 8).Value = 7625.25  # H64: 8834 -> 7625.25
$ws.Cells.Item(64, 9).Value = 7625.25  # I64: 6001.25 -> 7625.25
$ws.Cells.Item(64, 10).Value = 0  # J64: 14499.5 -> 0
$ws.Cells.Item(64, 11).Value = 22875.75  # K64: 18003.75 -> 22875.75
$ws.Cells.Item(64, 12).Value = 0  # L64: 43498.5 -> 0
$ws.Cells.Item(64, 13).Value = -22605.75  # M64: -17733.75 -> -22605.75
$ws.Cells.Item(64, 14).ClearContents()  # N64: -44038.5 -> (cleared)
$ws.Cells.Item(65, 8).Value = 5250  # H65: 0 -> 5250
$ws.Cells.Item(65, 9).Value = 3000  # I65: 0 -> 3000
$ws.Cells.Item(65, 10).Value = 7500  # J65: 0 -> 7500
$ws.Cells.Item(65, 11).Value = 27000  # K65: 0 -> 27000
$ws.Cells.Item(65, 12).Value = 67500  # L65: 0 -> 67500
$ws.Cells.Item(65, 13).Value = -23568  # M65: None -> -23568
$ws.Cells.Item(65, 14).Value = -74364  # N65: None -> -74364
$ws.Cells.Item(66, 8).Value = 3725  # H66: 0 -> 3725
$ws.Cells.Item(66, 9).Value = 4000  # I66: 0 -> 4000
$ws.Cells.Item(66, 10).Value = 3450  # J66: 0 -> 3450
$ws.Cells.Item(66, 11).Value = 36000  # K66: 0 -> 36000
$ws.Cells.Item(66, 12).Value = 31050  # L66: 0 -> 31050
$ws.Cells.Item(66, 13).Value = -32256  # M66: None -> -32256
$ws.Cells.Item(66, 14).Value = -38538  # N66: None -> -38538
$ws.Cells.Item(67, 8).Value = 7625.25  # H67: 8834 -> 7625.25
$ws.Cells.Item(67, 9).Value = 7625.25  # I67: 6001.25 -> 7625.25
$ws.Cells.Item(67, 10).Value = 0  # J67: 14499.5 -> 0
$ws.Cells.Item(67, 11).Value = 22875.75  # K67: 18003.75 -> 22875.75
$ws.Cells.Item(67, 12).Value = 0  # L67: 43498.5 -> 0
$ws.Cells.Item(67, 13).Value = -21939.75  # M67: -17067.75 -> -21939.75
$ws.Cells.Item(67, 14).ClearContents()  # N67: -45370.5 -> (cleared)
$ws.Cells.Item(68, 8).Value = 1248.1428  # H68: 1276.8572 -> 1248.1428
$ws.Cells.Item(68, 10).Value = 1481  # J68: 1538.4286 -> 1481
$ws.Cells.Item(68, 12).Value = 4443  # L68: 4615.2858 -> 4443
$ws.Cells.Item(68, 14).Value = -6065  # N68: -6237.2858 -> -6065
$ws.Cells.Item(70, 8).Value = 11740.333  # H70: 14812.546 -> 11740.333
$ws.Cells.Item(70, 9).Value = 7611  # I70: 6986 -> 7611
$ws.Cells.Item(70, 10).Value = 19999  # J70: 19284.857 -> 19999
$ws.Cells.Item(70, 11).Value = 22833  # K70: 20958 -> 22833
$ws.Cells.Item(70, 12).Value = 59997  # L70: 57854.571 -> 59997
$ws.Cells.Item(70, 13).Value = -22518  # M70: -20643 -> -22518
$ws.Cells.Item(70, 14).Value = -60627  # N70: -58484.571 -> -60627
$ws.Cells.Item(71, 8).Value = 1248.1428  # H71: 1276.8572 -> 1248.1428
$ws.Cells.Item(71, 10).Value = 1481  # J71: 1538.4286 -> 1481
$ws.Cells.Item(71, 12).Value = 13329  # L71: 13845.8574 -> 13329
$ws.Cells.Item(71, 14).Value = -21441  # N71: -21957.8574 -> -21441
$ws.Cells.Item(73, 8).Value = 11740.333  # H73: 14812.546 -> 11740.333
$ws.Cells.Item(73, 9).Value = 7611  # I73: 6986 -> 7611
$ws.Cells.Item(73, 10).Value = 19999  # J73: 19284.857 -> 19999
$ws.Cells.Item(73, 11).Value = 22833  # K73: 20958 -> 22833
$ws.Cells.Item(73, 12).Value = 59997  # L73: 57854.571 -> 59997
$ws.Cells.Item(73, 13).Value = -21741  # M73: -19866 -> -21741
$ws.Cells.Item(73, 14).Value = -62181  # N73: -60038.571 -> -62181
$ws.Cells.Item(76, 8).Value = 19487.5  # H76: 21784.143 -> 19487.5
$ws.Cells.Item(76, 9).Value = 19487.5  # I76: 19372.25 -> 19487.5
$ws.Cells.Item(76, 10).Value = 0  # J76: 25000 -> 0
$ws.Cells.Item(76, 11).Value = 58462.5  # K76: 58116.75 -> 58462.5
$ws.Cells.Item(76, 12).Value = 0  # L76: 75000 -> 0
$ws.Cells.Item(76, 13).Value = -58079.5  # M76: -57733.75 -> -58079.5
$ws.Cells.Item(76, 14).ClearContents()  # N76: -75766 -> (cleared)
$ws.Cells.Item(79, 8).Value = 19487.5  # H79: 21784.143 -> 19487.5
$ws.Cells.Item(79, 9).Value = 19487.5  # I79: 19372.25 -> 19487.5
$ws.Cells.Item(79, 10).Value = 0  # J79: 25000 -> 0
$ws.Cells.Item(79, 11).Value = 58462.5  # K79: 58116.75 -> 58462.5
$ws.Cells.Item(79, 12).Value = 0  # L79: 75000 -> 0
$ws.Cells.Item(79, 13).Value = -57136.5  # M79: -56790.75 -> -57136.5
$ws.Cells.Item(79, 14).ClearContents()  # N79: -77652 -> (cleared)
$ws.Cells.Item(80, 8).Value = 3743.75  # H80: 4622.5 -> 3743.75
$ws.Cells.Item(80, 10).Value = 3743.75  # J80: 4622.5 -> 3743.75
$ws.Cells.Item(80, 12).Value = 11231.25  # L80: 13867.5 -> 11231.25
$ws.Cells.Item(80, 14).Value = -13103.25  # N80: -15739.5 -> -13103.25
$ws.Cells.Item(81, 8).Value = 3750  # H81: 0 -> 3750
$ws.Cells.Item(81, 9).Value = 3000  # I81: 0 -> 3000
$ws.Cells.Item(81, 10).Value = 4500  # J81: 0 -> 4500
$ws.Cells.Item(81, 11).Value = 9000  # K81: 0 -> 9000
$ws.Cells.Item(81, 12).Value = 13500  # L81: 0 -> 13500
$ws.Cells.Item(81, 13).Value = -7877  # M81: None -> -7877
$ws.Cells.Item(81, 14).Value = -15746  # N81: None -> -15746
$ws.Cells.Item(83, 8).Value = 3743.75  # H83: 4622.5 -> 3743.75
$ws.Cells.Item(83, 10).Value = 3743.75  # J83: 4622.5 -> 3743.75
$ws.Cells.Item(83, 12).Value = 33693.75  # L83: 41602.5 -> 33693.75
$ws.Cells.Item(83, 14).Value = -43053.75  # N83: -50962.5 -> -43053.75
$ws.Cells.Item(84, 8).Value = 3750  # H84: 0 -> 3750
$ws.Cells.Item(84, 9).Value = 3000  # I84: 0 -> 3000
$ws.Cells.Item(84, 10).Value = 4500  # J84: 0 -> 4500
$ws.Cells.Item(84, 11).Value = 27000  # K84: 0 -> 27000
$ws.Cells.Item(84, 12).Value = 40500  # L84: 0 -> 40500
$ws.Cells.Item(84, 13).Value = -21384  # M84: None -> -21384
$ws.Cells.Item(84, 14).Value = -51732  # N84: None -> -51732
$ws.Cells.Item(86, 8).Value = 599  # H86: 598.6667 -> 599
$ws.Cells.Item(86, 9).Value = 0  # I86: 598 -> 0
$ws.Cells.Item(86, 11).Value = 0  # K86: 1794 -> 0
$ws.Cells.Item(86, 13).ClearContents()  # M86: -608 -> (cleared)
$ws.Cells.Item(89, 8).Value = 599  # H89: 598.6667 -> 599
$ws.Cells.Item(89, 9).Value = 0  # I89: 598 -> 0
$ws.Cells.Item(89, 11).Value = 0  # K89: 5382 -> 0
$ws.Cells.Item(89, 13).ClearContents()  # M89: 546 -> (cleared)
$ws.Cells.Item(98, 8).Value = 1257.7142  # H98: 1181.4546 -> 1257.7142
$ws.Cells.Item(98, 9).Value = 1015.8182  # I98: 999.9 -> 1015.8182
$ws.Cells.Item(98, 10).Value = 2144.6667  # J98: 2997 -> 2144.6667
$ws.Cells.Item(98, 11).Value = 3047.4546  # K98: 2999.7 -> 3047.4546
$ws.Cells.Item(98, 12).Value = 6434.000100000001  # L98: 8991 -> 6434.000100000001
$ws.Cells.Item(98, 13).Value = -1549.4546  # M98: -1501.7 -> -1549.4546
$ws.Cells.Item(98, 14).Value = -9430.000100000001  # N98: -11987 -> -9430.000100000001
$ws.Cells.Item(131, 8).Value = 1854.421  # H131: 1896.2354 -> 1854.421
$ws.Cells.Item(131, 10).Value = 2299.7273  # J131: 2477.6667 -> 2299.7273
$ws.Cells.Item(131, 12).Value = 6899.1819  # L131: 7433.000100000001 -> 6899.1819
$ws.Cells.Item(131, 14).Value = -16979.1819  # N131: -17513.0001 -> -16979.1819
$ws.Cells.Item(133, 8).Value = 9365.5  # H133: 9600.875 -> 9365.5
$ws.Cells.Item(133, 9).Value = 9365.5  # I133: 9600.875 -> 9365.5
$ws.Cells.Item(133, 11).Value = 28096.5  # K133: 28802.625 -> 28096.5
$ws.Cells.Item(133, 13).Value = -23036.5  # M133: -23742.625 -> -23036.5
$ws.Cells.Item(138, 8).Value = 6421.6665  # H138: 6600.1113 -> 6421.6665
$ws.Cells.Item(138, 9).Value = 2542.4285  # I138: 2771.8572 -> 2542.4285
$ws.Cells.Item(138, 11).Value = 7627.2855  # K138: 8315.5716 -> 7627.2855
$ws.Cells.Item(138, 13).Value = -2487.2855  # M138: -3175.571599999999 -> -2487.2855
$ws.Cells.Item(139, 8).Value = 5248.7144  # H139: 4901.2856 -> 5248.7144
$ws.Cells.Item(139, 9).Value = 3290.111  # I139: 3103.6 -> 3290.111
$ws.Cells.Item(139, 10).Value = 8774.2  # J139: 9395.5 -> 8774.2
$ws.Cells.Item(139, 11).Value = 9870.332999999999  # K139: 9310.8 -> 9870.332999999999
$ws.Cells.Item(139, 12).Value = 26322.6  # L139: 28186.5 -> 26322.6
$ws.Cells.Item(139, 13).Value = -4730.332999999999  # M139: -4170.799999999999 -> -4730.332999999999
$ws.Cells.Item(139, 14).Value = -36602.60000000001  # N139: -38466.5 -> -36602.60000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(26, 8).Value = 0  # H26: 22000 -> 0
$ws.Cells.Item(26, 10).Value = 0  # J26: 22000 -> 0
$ws.Cells.Item(26, 12).Value = 0  # L26: 22000 -> 0
$ws.Cells.Item(26, 14).ClearContents()  # N26: -22560 -> (cleared)
$ws.Cells.Item(50, 8).Value = 0  # H50: 22000 -> 0
$ws.Cells.Item(50, 10).Value = 0  # J50: 22000 -> 0
$ws.Cells.Item(50, 12).Value = 0  # L50: 22000 -> 0
$ws.Cells.Item(50, 14).ClearContents()  # N50: -22996 -> (cleared)
$ws.Cells.Item(104, 8).Value = 0  # H104: 69000 -> 0
$ws.Cells.Item(104, 10).Value = 0  # J104: 69000 -> 0
$ws.Cells.Item(104, 12).Value = 0  # L104: 69000 -> 0
$ws.Cells.Item(104, 14).ClearContents()  # N104: -75988 -> (cleared)
$ws.Cells.Item(127, 8).Value = 39990  # H127: 56244.5 -> 39990
$ws.Cells.Item(127, 10).Value = 39990  # J127: 56244.5 -> 39990
$ws.Cells.Item(127, 12).Value = 39990  # L127: 56244.5 -> 39990
$ws.Cells.Item(127, 14).Value = -49910  # N127: -66164.5 -> -49910
$ws.Cells.Item(132, 8).Value = 1812.3208  # H132: 1798.8889 -> 1812.3208
$ws.Cells.Item(132, 9).Value = 1772.3062  # I132: 1776.4082 -> 1772.3062
$ws.Cells.Item(132, 10).Value = 2302.5  # J132: 2019.2 -> 2302.5
$ws.Cells.Item(132, 11).Value = 5316.9186  # K132: 5329.2246 -> 5316.9186
$ws.Cells.Item(132, 12).Value = 6907.5  # L132: 6057.6 -> 6907.5
$ws.Cells.Item(132, 13).Value = -2786.9186  # M132: -2799.2246 -> -2786.9186
$ws.Cells.Item(132, 14).Value = -11967.5  # N132: -11117.6 -> -11967.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value = 3000  # H22: 2679.5 -> 3000
$ws.Cells.Item(22, 9).Value = 0  # I22: 2359 -> 0
$ws.Cells.Item(22, 11).Value = 0  # K22: 2359 -> 0
$ws.Cells.Item(22, 13).ClearContents()  # M22: -2064 -> (cleared)
$ws.Cells.Item(27, 8).Value = 3000  # H27: 2679.5 -> 3000
$ws.Cells.Item(27, 9).Value = 0  # I27: 2359 -> 0
$ws.Cells.Item(27, 11).Value = 0  # K27: 2359 -> 0
$ws.Cells.Item(27, 13).ClearContents()  # M27: -2252 -> (cleared)
$ws.Cells.Item(46, 8).Value = 6987.1113  # H46: 5648.5835 -> 6987.1113
$ws.Cells.Item(46, 9).Value = 6987.1113  # I46: 6438.4 -> 6987.1113
$ws.Cells.Item(46, 10).Value = 0  # J46: 1699.5 -> 0
$ws.Cells.Item(46, 11).Value = 6987.1113  # K46: 6438.4 -> 6987.1113
$ws.Cells.Item(46, 12).Value = 0  # L46: 1699.5 -> 0
$ws.Cells.Item(46, 13).Value = -6799.1113  # M46: -6250.4 -> -6799.1113
$ws.Cells.Item(46, 14).ClearContents()  # N46: -2075.5 -> (cleared)
$ws.Cells.Item(55, 8).Value = 362.41177  # H55: 448.0909 -> 362.41177
$ws.Cells.Item(55, 9).Value = 189.77777  # I55: 204.83333 -> 189.77777
$ws.Cells.Item(55, 10).Value = 556.625  # J55: 740 -> 556.625
$ws.Cells.Item(55, 11).Value = 189.77777  # K55: 204.83333 -> 189.77777
$ws.Cells.Item(55, 12).Value = 556.625  # L55: 740 -> 556.625
$ws.Cells.Item(55, 13).Value = -16.77777  # M55: -31.83332999999999 -> -16.77777
$ws.Cells.Item(55, 14).Value = -902.625  # N55: -1086 -> -902.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(37, 8).Value = 2992.5  # H37: 2995 -> 2992.5
$ws.Cells.Item(37, 9).Value = 2990  # I37: 0 -> 2990
$ws.Cells.Item(37, 11).Value = 2990  # K37: 0 -> 2990
$ws.Cells.Item(37, 13).Value = -2787  # M37: None -> -2787
$ws.Cells.Item(62, 8).Value = 11371.857  # H62: 8273.77 -> 11371.857
$ws.Cells.Item(62, 9).Value = 9937.5  # I62: 8104.6665 -> 9937.5
$ws.Cells.Item(62, 10).Value = 13284.333  # J62: 8418.714 -> 13284.333
$ws.Cells.Item(62, 11).Value = 9937.5  # K62: 8104.6665 -> 9937.5
$ws.Cells.Item(62, 12).Value = 13284.333  # L62: 8418.714 -> 13284.333
$ws.Cells.Item(62, 13).Value = -9313.5  # M62: -7480.6665 -> -9313.5
$ws.Cells.Item(62, 14).Value = -14532.333  # N62: -9666.714 -> -14532.333
$ws.Cells.Item(65, 8).Value = 11371.857  # H65: 8273.77 -> 11371.857
$ws.Cells.Item(65, 9).Value = 9937.5  # I65: 8104.6665 -> 9937.5
$ws.Cells.Item(65, 10).Value = 13284.333  # J65: 8418.714 -> 13284.333
$ws.Cells.Item(65, 11).Value = 49687.5  # K65: 40523.3325 -> 49687.5
$ws.Cells.Item(65, 12).Value = 66421.66500000001  # L65: 42093.57 -> 66421.66500000001
$ws.Cells.Item(65, 13).Value = -46567.5  # M65: -37403.3325 -> -46567.5
$ws.Cells.Item(65, 14).Value = -72661.66500000001  # N65: -48333.57 -> -72661.66500000001
$ws.Cells.Item(81, 8).Value = 3549.1428  # H81: 2684.2856 -> 3549.1428
$ws.Cells.Item(81, 9).Value = 3549.1428  # I81: 2724 -> 3549.1428
$ws.Cells.Item(81, 10).Value = 0  # J81: 2446 -> 0
$ws.Cells.Item(81, 11).Value = 7098.2856  # K81: 5448 -> 7098.2856
$ws.Cells.Item(81, 12).Value = 0  # L81: 4892 -> 0
$ws.Cells.Item(81, 13).Value = -6037.2856  # M81: -4387 -> -6037.2856
$ws.Cells.Item(81, 14).ClearContents()  # N81: -7014 -> (cleared)
$ws.Cells.Item(84, 8).Value = 3549.1428  # H84: 2684.2856 -> 3549.1428
$ws.Cells.Item(84, 9).Value = 3549.1428  # I84: 2724 -> 3549.1428
$ws.Cells.Item(84, 10).Value = 0  # J84: 2446 -> 0
$ws.Cells.Item(84, 11).Value = 35491.428  # K84: 27240 -> 35491.428
$ws.Cells.Item(84, 12).Value = 0  # L84: 24460 -> 0
$ws.Cells.Item(84, 13).Value = -30187.428  # M84: -21936 -> -30187.428
$ws.Cells.Item(84, 14).ClearContents()  # N84: -35068 -> (cleared)
$ws.Cells.Item(117, 8).Value = 0  # H117: 110200 -> 0
$ws.Cells.Item(117, 10).Value = 0  # J117: 110200 -> 0
$ws.Cells.Item(117, 12).Value = 0  # L117: 110200 -> 0
$ws.Cells.Item(117, 14).ClearContents()  # N117: -119378 -> (cleared)
$ws.Cells.Item(139, 8).Value = 104524.164  # H139: 101163.57 -> 104524.164
